$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header formatting (bold, border, alignment) from H1 so the new
# header cells match the existing header style exactly.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-39: add I (I0) and J (IF) columns ---
# I is 1 for every row except row 38, which is 4.
# J mirrors the existing H (IP) value for every row except row 38,
# where it is 6 instead of H38's 3.
$iVals = @{ 38 = 4 }
$jOverrides = @{ 38 = 6 }

for ($r = 2; $r -le 39; $r++) {
    if ($iVals.ContainsKey($r)) {
        $iValue = $iVals[$r]
    } else {
        $iValue = 1
    }
    $ws.Cells.Item($r, 9).Value = $iValue

    if ($jOverrides.ContainsKey($r)) {
        $jValue = $jOverrides[$r]
    } else {
        $jValue = $ws.Cells.Item($r, 8).Value()
    }
    $ws.Cells.Item($r, 10).Value = $jValue
}
